# Auto-generated Excel COM-interop script
# Applies row-level odds data corrections to the "Colombia Primera B" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 176
$ws.Range('B176').Value2 = 7630528
$ws.Range('F176').Value2 = 'Atletico Cali FC'
$ws.Range('G176').Value2 = 'Cucuta Deportivo'
$ws.Range('H176').Value2 = 0
$ws.Range('J176').Value2 = 'A'
$ws.Range('K176').Value2 = 3.2
$ws.Range('L176').Value2 = 3.1
$ws.Range('M176').Value2 = 2.2
$ws.Range('N176').Value2 = 3.3
$ws.Range('O176').Value2 = 3.25
$ws.Range('P176').Value2 = 2.2
$ws.Range('Q176').Value2 = 0.25
$ws.Range('R176').Value2 = 1.9
$ws.Range('S176').Value2 = 1.9
$ws.Range('U176').Value2 = 1.775
$ws.Range('V176').Value2 = 2.025
$ws.Range('X176').Value2 = -1
$ws.Range('Y176').Value2 = 1.2
$ws.Range('Z176').Value2 = -1
$ws.Range('AA176').Value2 = 0.8999999999999999
$ws.Range('AB176').Value2 = -1
$ws.Range('AC176').Value2 = 1.025

# Row 177
$ws.Range('B177').Value2 = 7630522
$ws.Range('F177').Value2 = 'Orsomarso'
$ws.Range('G177').Value2 = 'Real Soacha Cundinamarca'
$ws.Range('H177').Value2 = 1
$ws.Range('J177').Value2 = 'D'
$ws.Range('K177').Value2 = 2.2
$ws.Range('L177').Value2 = 3
$ws.Range('M177').Value2 = 3.2
$ws.Range('N177').Value2 = 2.375
$ws.Range('O177').Value2 = 3.2
$ws.Range('P177').Value2 = 3.2
$ws.Range('Q177').Value2 = -0.25
$ws.Range('R177').Value2 = 2.025
$ws.Range('S177').Value2 = 1.775
$ws.Range('U177').Value2 = 1.975
$ws.Range('V177').Value2 = 1.825
$ws.Range('X177').Value2 = 2.2
$ws.Range('Y177').Value2 = -1
$ws.Range('Z177').Value2 = -0.5
$ws.Range('AA177').Value2 = 0.3875
$ws.Range('AB177').Value2 = 0
$ws.Range('AC177').Value2 = -0

# Row 188
$ws.Range('B188').Value2 = 7630542
$ws.Range('F188').Value2 = 'Orsomarso'
$ws.Range('G188').Value2 = 'Union Magdalena'
$ws.Range('H188').Value2 = 1
$ws.Range('I188').Value2 = 3
$ws.Range('J188').Value2 = 'A'
$ws.Range('K188').Value2 = 2.5
$ws.Range('L188').Value2 = 3.2
$ws.Range('M188').Value2 = 2.5
$ws.Range('N188').Value2 = 3
$ws.Range('O188').Value2 = 3.1
$ws.Range('P188').Value2 = 2.15
$ws.Range('Q188').Value2 = 0.25
$ws.Range('R188').Value2 = 1.85
$ws.Range('S188').Value2 = 1.95
$ws.Range('T188').Value2 = 2
$ws.Range('U188').Value2 = 1.75
$ws.Range('V188').Value2 = 2.05
$ws.Range('X188').Value2 = -1
$ws.Range('Y188').Value2 = 1.15
$ws.Range('Z188').Value2 = -1
$ws.Range('AA188').Value2 = 0.95
$ws.Range('AB188').Value2 = 0.75

# Row 189
$ws.Range('B189').Value2 = 7630540
$ws.Range('F189').Value2 = 'Barranquilla FC'
$ws.Range('G189').Value2 = 'Leones'
$ws.Range('H189').Value2 = 2
$ws.Range('I189').Value2 = 2
$ws.Range('J189').Value2 = 'D'
$ws.Range('K189').Value2 = 3.6
$ws.Range('L189').Value2 = 3.5
$ws.Range('M189').Value2 = 2
$ws.Range('N189').Value2 = 2.1
$ws.Range('O189').Value2 = 3.5
$ws.Range('P189').Value2 = 3.4
$ws.Range('Q189').Value2 = -0.25
$ws.Range('R189').Value2 = 1.825
$ws.Range('S189').Value2 = 1.975
$ws.Range('T189').Value2 = 2.5
$ws.Range('U189').Value2 = 1.8
$ws.Range('V189').Value2 = 2
$ws.Range('X189').Value2 = 2.5
$ws.Range('Y189').Value2 = -1
$ws.Range('Z189').Value2 = -0.5
$ws.Range('AA189').Value2 = 0.4875
$ws.Range('AB189').Value2 = 0.8

# Row 227
$ws.Range('B227').Value2 = 7630576
$ws.Range('F227').Value2 = 'Real San Andres'
$ws.Range('G227').Value2 = 'Union Magdalena'
$ws.Range('H227').Value2 = 2
$ws.Range('I227').Value2 = 0
$ws.Range('J227').Value2 = 'H'
$ws.Range('K227').Value2 = 2.75
$ws.Range('L227').Value2 = 3.1
$ws.Range('M227').Value2 = 2.4
$ws.Range('N227').Value2 = 3.4
$ws.Range('O227').Value2 = 3.6
$ws.Range('P227').Value2 = 2.05
$ws.Range('Q227').Value2 = 0.5
$ws.Range('R227').Value2 = 1.775
$ws.Range('S227').Value2 = 2.025
$ws.Range('T227').Value2 = 2.75
$ws.Range('U227').Value2 = 1.95
$ws.Range('V227').Value2 = 1.85
$ws.Range('W227').Value2 = 2.4
$ws.Range('Y227').Value2 = -1
$ws.Range('Z227').Value2 = 0.7749999999999999
$ws.Range('AA227').Value2 = -1
$ws.Range('AB227').Value2 = -1
$ws.Range('AC227').Value2 = 0.8500000000000001

# Row 228
$ws.Range('B228').Value2 = 7630580
$ws.Range('F228').Value2 = 'Atletico Cali FC'
$ws.Range('G228').Value2 = 'Orsomarso'
$ws.Range('H228').Value2 = 0
$ws.Range('I228').Value2 = 2
$ws.Range('J228').Value2 = 'A'
$ws.Range('K228').Value2 = 2
$ws.Range('L228').Value2 = 3.25
$ws.Range('M228').Value2 = 3.6
$ws.Range('N228').Value2 = 2.4
$ws.Range('O228').Value2 = 3.1
$ws.Range('P228').Value2 = 3.2
$ws.Range('Q228').Value2 = -0.25
$ws.Range('R228').Value2 = 2
$ws.Range('S228').Value2 = 1.8
$ws.Range('T228').Value2 = 2
$ws.Range('U228').Value2 = 1.85
$ws.Range('V228').Value2 = 1.95
$ws.Range('W228').Value2 = -1
$ws.Range('Y228').Value2 = 2.2
$ws.Range('Z228').Value2 = -1
$ws.Range('AA228').Value2 = 0.8
$ws.Range('AB228').Value2 = 0
$ws.Range('AC228').Value2 = -0

# Row 243
$ws.Range('B243').Value2 = 7658164
$ws.Range('E243').Value2 = 45385.90277777778
$ws.Range('F243').Value2 = 'Leones'
$ws.Range('G243').Value2 = 'Orsomarso'
$ws.Range('K243').Value2 = 2.3
$ws.Range('M243').Value2 = 3
$ws.Range('N243').Value2 = 2.3
$ws.Range('O243').Value2 = 3.4
$ws.Range('P243').Value2 = 3.2
$ws.Range('R243').Value2 = 1.975
$ws.Range('S243').Value2 = 1.825
$ws.Range('T243').Value2 = 2.5
$ws.Range('U243').Value2 = 1.925
$ws.Range('V243').Value2 = 1.875

# Row 244
$ws.Range('B244').Value2 = 7658166
$ws.Range('E244').Value2 = 45386.72916666666
$ws.Range('F244').Value2 = 'Atletico Cali FC'
$ws.Range('G244').Value2 = 'Real Cartagena'
$ws.Range('K244').Value2 = 3.1
$ws.Range('L244').Value2 = 3
$ws.Range('M244').Value2 = 2.4
$ws.Range('N244').Value2 = 4.5
$ws.Range('O244').Value2 = 3.25
$ws.Range('P244').Value2 = 1.909
$ws.Range('Q244').Value2 = 0.5
$ws.Range('R244').Value2 = 1.925
$ws.Range('S244').Value2 = 1.875
$ws.Range('T244').Value2 = 2.25
$ws.Range('U244').Value2 = 1.875
$ws.Range('V244').Value2 = 1.925

# Row 245
$ws.Range('B245').Value2 = 7657924
$ws.Range('E245').Value2 = 45386.83333333334
$ws.Range('F245').Value2 = 'Real Soacha Cundinamarca'
$ws.Range('G245').Value2 = 'Union Magdalena'
$ws.Range('K245').Value2 = 1.95
$ws.Range('L245').Value2 = 3.25
$ws.Range('M245').Value2 = 3.8
$ws.Range('N245').Value2 = 1.95
$ws.Range('O245').Value2 = 3.2
$ws.Range('P245').Value2 = 3.75
$ws.Range('Q245').Value2 = -0.5
$ws.Range('R245').Value2 = 1.95
$ws.Range('S245').Value2 = 1.85
$ws.Range('T245').Value2 = 2.5
$ws.Range('U245').Value2 = 1.875
$ws.Range('V245').Value2 = 1.925

# Row 246
$ws.Range('R246').Value2 = 2
$ws.Range('S246').Value2 = 1.8

# Row 247
$ws.Range('N247').Value2 = 1.85
$ws.Range('O247').Value2 = 3.25
$ws.Range('P247').Value2 = 4
$ws.Range('R247').Value2 = 1.875
$ws.Range('S247').Value2 = 1.925
